$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 88: Coffeeland ---
$ws.Range("A88").Value = "Coffeeland"
$ws.Range("B88").Value = "Augustine Sedgewick"
$ws.Range("C87").Copy()
$ws.Range("C88").PasteSpecial(-4122)
$ws.Range("C88").Value = 43992
$ws.Range("D87").Copy()
$ws.Range("D88").PasteSpecial(-4122)
$ws.Range("D88").Value = 43994
$ws.Range("E88").Value = "coffee;business;socialism;capitalism;exploitation;history"
$ws.Range("F88").Value = "Audio"
$ws.Range("G88").Value = "15 Hours 2 Mins"

# --- Row 89: Brick by Brick ---
$ws.Range("A89").Value = "Brick by Brick"
$ws.Range("B89").Value = "David Robertson"
$ws.Range("C89").Value = "6/12/20202"
$ws.Range("D87").Copy()
$ws.Range("D89").PasteSpecial(-4122)
$ws.Range("D89").Value = 43996
$ws.Range("E89").Value = "Lego;business;history;innovation;toys"
$ws.Range("F89").Value = "Audio"
$ws.Range("G89").Value = "10 Hours 24 Mins"

# --- Row 90: The Virgin Way ---
$ws.Range("A90").Value = "The Virgin Way"
$ws.Range("B90").Value = "Richard Branson"
$ws.Range("C87").Copy()
$ws.Range("C90").PasteSpecial(-4122)
$ws.Range("C90").Value = 43997
$ws.Range("D87").Copy()
$ws.Range("D90").PasteSpecial(-4122)
$ws.Range("D90").Value = 43998
$ws.Range("E90").Value = "richard branson;business;leadership;virgin;airlines"
$ws.Range("F90").Value = "Audio"
$ws.Range("G90").Value = "11 Hours 30 Mins"

$excel.CutCopyMode = $false

# Update the visible selection/active cell to mirror the saved view state
$null = $ws.Range("A91").Select()
